$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spring 2022 course list (column C) shifts due to a newly added course:
# C4: GEOL 3225 -> CYBR 4125
# C5: CYBR 4125 -> CPSC 4135
# C6: CPSC 4135 -> CYBR 4145 (new course inserted)
# C7 / C8 keep their existing values (CPSC 4148 / CPSC 4155)
$ws.Range("C4").Value = "CYBR 4125"
$ws.Range("C5").Value = "CPSC 4135"
$ws.Range("C6").Value = "CYBR 4145"

# Fall 2022 course list (columns A/B):
# A5: PSYC 1105 -> GEOL 1121, credits 2 -> 3
$ws.Range("A5").Value = "GEOL 1121"
$ws.Range("B5").Value = 3

# A9: CPSC 3415 -> CPSC 4000, credits 1 -> 0
$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0

# Row 10 (old CPSC 4000 / 0 credits) is removed entirely
$ws.Range("A10:B10").ClearContents()
